$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = @{ B = 0.2380952380952381; C = 0.5238095238095238; J = 0.04761904761904762; P = 0.09523809523809523; S = 0.09523809523809523 }
    3  = @{ P = 0.6363636363636364; S = 0.3636363636363636 }
    4  = @{ P = 0.5; S = 0.5 }
    6  = @{ B = 0.1666666666666667; J = 0.1666666666666667; O = 0.1666666666666667; S = 0.5 }
    7  = @{ J = 0.5; Q = 0.25; R = 0.25 }
    8  = @{ B = 0.04347826086956522; D = 0.08695652173913043; F = 0.08695652173913043; J = 0.08695652173913043; Q = 0.2173913043478261; S = 0.4782608695652174 }
    9  = @{ B = 0.25; F = 0.125; J = 0.125; Q = 0.125; S = 0.375 }
    10 = @{ B = 0.2352941176470588; D = 0.0392156862745098; F = 0.0196078431372549; J = 0.1176470588235294; O = 0.0196078431372549; Q = 0.2352941176470588; R = 0.0196078431372549; S = 0.3137254901960784 }
    11 = @{ G = 0.07692307692307693; J = 0.2307692307692308; K = 0.3846153846153846; L = 0.3076923076923077 }
    12 = @{ G = 0.75; J = 0.25 }
    15 = @{ H = 0.2; J = 0.3; S = 0.5 }
    16 = @{ I = 0.09090909090909091; J = 0.4545454545454545; K = 0.1818181818181818; S = 0.2727272727272727 }
    17 = @{ H = 0.4210526315789473; J = 0.3157894736842105; O = 0.1052631578947368; S = 0.1578947368421053 }
    18 = @{ H = 0.5; J = 0.5 }
    19 = @{ F = 0.01612903225806452; H = 0.1935483870967742; I = 0.1129032258064516; J = 0.3225806451612903; K = 0.1129032258064516; O = 0.08064516129032258; S = 0.1612903225806452 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $addr = "$col$row"
        $ws.Range($addr).Value = $cols[$col]
    }
}
